$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (233) down into the
# new rows (234:238) so the date column keeps its date style, then fill
# in the actual values for the continued daily series.
$ws.Range("A233:D233").Copy() | Out-Null
$ws.Range("A234:D238").PasteSpecial(-4122) | Out-Null

$newRows = @(
    @(44308, 0, 1, 48.07692307692308),
    @(44309, 0, 1, 48.07692307692308),
    @(44310, 0, 1, 48.07692307692308),
    @(44311, 0, 1, 48.07692307692308),
    @(44312, 0, 0, 0)
)

$startRow = 234
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
}
